$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3309.875
$ws.Range("I33").Value = 4332.3335
$ws.Range("K33").Value = 4332.3335
$ws.Range("M33").Value = -4103.3335
$ws.Range("H62").Value = 2894.4546
$ws.Range("I62").Value = 2486.4783
$ws.Range("J62").Value = 3832.8
$ws.Range("K62").Value = 2486.4783
$ws.Range("L62").Value = 3832.8
$ws.Range("M62").Value = -1862.4783
$ws.Range("N62").Value = -5080.8
$ws.Range("H65").Value = 2894.4546
$ws.Range("I65").Value = 2486.4783
$ws.Range("J65").Value = 3832.8
$ws.Range("K65").Value = 12432.3915
$ws.Range("L65").Value = 19164
$ws.Range("M65").Value = -9312.391500000002
$ws.Range("N65").Value = -25404
$ws.Range("H74").Value = 3353.25
$ws.Range("I74").Value = 3191.6875
$ws.Range("J74").Value = 3999.5
$ws.Range("K74").Value = 3191.6875
$ws.Range("L74").Value = 3999.5
$ws.Range("M74").Value = -2255.6875
$ws.Range("N74").Value = -5871.5
$ws.Range("H77").Value = 3353.25
$ws.Range("I77").Value = 3191.6875
$ws.Range("J77").Value = 3999.5
$ws.Range("K77").Value = 15958.4375
$ws.Range("L77").Value = 19997.5
$ws.Range("M77").Value = -11278.4375
$ws.Range("N77").Value = -29357.5
$ws.Range("H96").Value = 31252304
$ws.Range("J96").Value = 1907.25
$ws.Range("L96").Value = 5721.75
$ws.Range("N96").Value = -8467.75
$ws.Range("H120").Value = 34187.5
$ws.Range("J120").Value = 34187.5
$ws.Range("L120").Value = 34187.5
$ws.Range("N120").Value = -43863.5
$ws.Range("H132").Value = 1397.3715
$ws.Range("I132").Value = 1360.5454
$ws.Range("K132").Value = 4081.6362
$ws.Range("M132").Value = -1551.6362
$ws.Range("H135").Value = 510.65518
$ws.Range("I135").Value = 457.4643
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 4117.178699999999
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -1582.178699999999
$ws.Range("N135").Value = -23070
$ws.Range("H141").Value = 8286.541
$ws.Range("I141").Value = 1243.6086
$ws.Range("J141").Value = 19857.072
$ws.Range("K141").Value = 3730.8258
$ws.Range("L141").Value = 59571.216
$ws.Range("M141").Value = 1449.1742
$ws.Range("N141").Value = -69931.216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1643.65
$ws.Range("I2").Value = 1769
$ws.Range("J2").Value = 933.3333
$ws.Range("K2").Value = 1769
$ws.Range("L2").Value = 933.3333
$ws.Range("M2").Value = -1656
$ws.Range("N2").Value = -1159.3333
$ws.Range("H32").Value = 7378.1274
$ws.Range("I32").Value = 5531.0244
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 5531.0244
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -5244.0244
$ws.Range("N32").Value = -20574
$ws.Range("H45").Value = 1780.2142
$ws.Range("I45").Value = 1560.4286
$ws.Range("K45").Value = 1560.4286
$ws.Range("M45").Value = -1183.4286
$ws.Range("H74").Value = 1769.25
$ws.Range("I74").Value = 1741.0476
$ws.Range("J74").Value = 1966.6666
$ws.Range("K74").Value = 1741.0476
$ws.Range("L74").Value = 1966.6666
$ws.Range("M74").Value = -867.0476000000001
$ws.Range("N74").Value = -3714.6666
$ws.Range("H77").Value = 1769.25
$ws.Range("I77").Value = 1741.0476
$ws.Range("J77").Value = 1966.6666
$ws.Range("K77").Value = 8705.238000000001
$ws.Range("L77").Value = 9833.333
$ws.Range("M77").Value = -4337.238000000001
$ws.Range("N77").Value = -18569.333
$ws.Range("H116").Value = 1643.65
$ws.Range("I116").Value = 1769
$ws.Range("J116").Value = 933.3333
$ws.Range("K116").Value = 1769
$ws.Range("L116").Value = 933.3333
$ws.Range("M116").Value = 525
$ws.Range("N116").Value = -5521.3333
$ws.Range("H121").ClearContents()
$ws.Range("I121").ClearContents()
$ws.Range("J121").ClearContents()
$ws.Range("K121").ClearContents()
$ws.Range("L121").ClearContents()
$ws.Range("N121").ClearContents()
$ws.Range("H122").ClearContents()
$ws.Range("I122").ClearContents()
$ws.Range("J122").ClearContents()
$ws.Range("K122").ClearContents()
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H123").ClearContents()
$ws.Range("I123").ClearContents()
$ws.Range("J123").ClearContents()
$ws.Range("K123").ClearContents()
$ws.Range("L123").ClearContents()
$ws.Range("H124").ClearContents()
$ws.Range("I124").ClearContents()
$ws.Range("J124").ClearContents()
$ws.Range("K124").ClearContents()
$ws.Range("L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125").ClearContents()
$ws.Range("I125").ClearContents()
$ws.Range("J125").ClearContents()
$ws.Range("K125").ClearContents()
$ws.Range("L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H126").ClearContents()
$ws.Range("I126").ClearContents()
$ws.Range("J126").ClearContents()
$ws.Range("K126").ClearContents()
$ws.Range("L126").ClearContents()
$ws.Range("H127").ClearContents()
$ws.Range("I127").ClearContents()
$ws.Range("J127").ClearContents()
$ws.Range("K127").ClearContents()
$ws.Range("L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128").ClearContents()
$ws.Range("I128").ClearContents()
$ws.Range("J128").ClearContents()
$ws.Range("K128").ClearContents()
$ws.Range("L128").ClearContents()
$ws.Range("N128").ClearContents()
$ws.Range("H129").ClearContents()
$ws.Range("I129").ClearContents()
$ws.Range("J129").ClearContents()
$ws.Range("K129").ClearContents()
$ws.Range("L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130").ClearContents()
$ws.Range("I130").ClearContents()
$ws.Range("J130").ClearContents()
$ws.Range("K130").ClearContents()
$ws.Range("L130").ClearContents()
$ws.Range("H131").ClearContents()
$ws.Range("I131").ClearContents()
$ws.Range("J131").ClearContents()
$ws.Range("K131").ClearContents()
$ws.Range("L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132").ClearContents()
$ws.Range("I132").ClearContents()
$ws.Range("J132").ClearContents()
$ws.Range("K132").ClearContents()
$ws.Range("L132").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133").ClearContents()
$ws.Range("I133").ClearContents()
$ws.Range("J133").ClearContents()
$ws.Range("K133").ClearContents()
$ws.Range("L133").ClearContents()
$ws.Range("M133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134").ClearContents()
$ws.Range("I134").ClearContents()
$ws.Range("J134").ClearContents()
$ws.Range("K134").ClearContents()
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135").ClearContents()
$ws.Range("I135").ClearContents()
$ws.Range("J135").ClearContents()
$ws.Range("K135").ClearContents()
$ws.Range("L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H137").ClearContents()
$ws.Range("I137").ClearContents()
$ws.Range("J137").ClearContents()
$ws.Range("K137").ClearContents()
$ws.Range("L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138").ClearContents()
$ws.Range("I138").ClearContents()
$ws.Range("J138").ClearContents()
$ws.Range("K138").ClearContents()
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139").ClearContents()
$ws.Range("I139").ClearContents()
$ws.Range("J139").ClearContents()
$ws.Range("K139").ClearContents()
$ws.Range("L139").ClearContents()
$ws.Range("H140").ClearContents()
$ws.Range("I140").ClearContents()
$ws.Range("J140").ClearContents()
$ws.Range("K140").ClearContents()
$ws.Range("L140").ClearContents()
$ws.Range("H141").ClearContents()
$ws.Range("I141").ClearContents()
$ws.Range("J141").ClearContents()
$ws.Range("K141").ClearContents()
$ws.Range("L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1643.65
$ws.Range("I3").Value = 1769
$ws.Range("J3").Value = 933.3333
$ws.Range("K3").Value = 1769
$ws.Range("L3").Value = 933.3333
$ws.Range("M3").Value = -1655
$ws.Range("N3").Value = -1161.3333
$ws.Range("H6").Value = 19141.334
$ws.Range("J6").Value = 19141.334
$ws.Range("L6").Value = 19141.334
$ws.Range("N6").Value = -19367.334
$ws.Range("H94").Value = 771.6667
$ws.Range("I94").Value = 412.5
$ws.Range("J94").Value = 1490
$ws.Range("K94").Value = 412.5
$ws.Range("L94").Value = 1490
$ws.Range("M94").Value = 38.5
$ws.Range("N94").Value = -2392

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 636
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 810.4
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 810.4
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = -1510.4
$ws.Range("H86").Value = 3358.125
$ws.Range("I86").Value = 2710.6
$ws.Range("K86").Value = 2710.6
$ws.Range("M86").Value = -1587.6
$ws.Range("H89").Value = 3358.125
$ws.Range("I89").Value = 2710.6
$ws.Range("K89").Value = 13553
$ws.Range("M89").Value = -7937
$ws.Range("H107").Value = 694.4545
$ws.Range("I107").Value = 229.42857
$ws.Range("J107").Value = 911.4667
$ws.Range("K107").Value = 229.42857
$ws.Range("L107").Value = 911.4667
$ws.Range("M107").Value = 1690.57143
$ws.Range("N107").Value = -4751.4667
$ws.Range("H129").Value = 45899.7
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 45899.7
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 45899.7
$ws.Range("N129").Value = -55899.7
$ws.Range("H130").Value = 0
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("H131").Value = 32000
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 32000
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 32000
$ws.Range("N131").Value = -42080
$ws.Range("H132").Value = 2342.0833
$ws.Range("I132").Value = 957.5417
$ws.Range("J132").Value = 5111.1665
$ws.Range("K132").Value = 2872.6251
$ws.Range("L132").Value = 15333.4995
$ws.Range("M132").Value = -342.6251000000002
$ws.Range("N132").Value = -20393.4995
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 8109
$ws.Range("I134").Value = 8109
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 24327
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -21792
$ws.Range("H135").Value = 42836.93
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 42836.93
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 42836.93
$ws.Range("N135").Value = -52976.93
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 56865
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 56865
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 56865
$ws.Range("N138").Value = -67145
$ws.Range("H139").Value = 41033.332
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 41033.332
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 41033.332
$ws.Range("N139").Value = -51313.332
$ws.Range("H140").Value = 54173.453
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 54173.453
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 54173.453
$ws.Range("N140").Value = -64533.453
$ws.Range("H141").Value = 33519.23
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 33519.23
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 33519.23
$ws.Range("N141").Value = -43879.23

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 4249.032
$ws.Range("I131").Value = 407.69232
$ws.Range("J131").Value = 7023.3335
$ws.Range("K131").Value = 1223.07696
$ws.Range("L131").Value = 21070.0005
$ws.Range("M131").Value = 3816.92304
$ws.Range("N131").Value = -31150.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2003.4615
$ws.Range("I82").Value = 935.5
$ws.Range("J82").Value = 3249.4167
$ws.Range("K82").Value = 935.5
$ws.Range("L82").Value = 3249.4167
$ws.Range("M82").Value = -574.5
$ws.Range("N82").Value = -3971.4167
$ws.Range("H85").Value = 2003.4615
$ws.Range("I85").Value = 935.5
$ws.Range("J85").Value = 3249.4167
$ws.Range("K85").Value = 935.5
$ws.Range("L85").Value = 3249.4167
$ws.Range("M85").Value = 312.5
$ws.Range("N85").Value = -5745.4167
$ws.Range("H93").Value = 6912.263
$ws.Range("I93").Value = 12514.333
$ws.Range("K93").Value = 12514.333
$ws.Range("M93").Value = -11266.333
$ws.Range("H122").Value = 27779624
$ws.Range("I122").Value = 111111110
$ws.Range("J122").Value = 2461.6667
$ws.Range("K122").Value = 333333330
$ws.Range("L122").Value = 7385.000100000001
$ws.Range("M122").Value = -333330880
$ws.Range("N122").Value = -12285.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5466280.5
$ws.Range("I136").Value = 7247265.5
$ws.Range("J136").Value = 4593.6665
$ws.Range("K136").Value = 21741796.5
$ws.Range("L136").Value = 13780.9995
$ws.Range("M136").Value = -21739246.5
$ws.Range("N136").Value = -18880.9995
